# The trailing row of this export held a stray attribution/footer note
# ("data generated with OnlineDataGenerator ...") in column A only, while
# every real data row has values in columns A:E. That lone free-text row
# breaks a clean, typed load into Postgres, so it gets selected and wiped,
# leaving row 102 present but blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

[void]$ws.Select()
[void]$ws.Rows.Item(102).Select()
$ws.Cells.Item(102, 1).ClearContents()
